$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82 (shifts existing rows 82-150 down to 83-151)
$ws.Rows.Item(82).Insert()

# Row 82 is a "section header" style row like row 16 / row 109 - column A stays empty
$ws.Range("A82").Clear()

# B82: new dose-form term
$ws.Range("B82").Value = "oral/rectal suspension"

# C82: definition text with a plain lead-in sentence followed by an italicized
# usage note (mirrors the rich-text pattern used elsewhere in this column,
# e.g. the "oral/rectal solution" row).
$plainPart = "Liquid preparation consisting of a suspension intended for oral or rectal use. "
$italicPart = "This term is only to be used in cases where there is not a single predominant route of administration for the medicinal product."
$ws.Range("C82").Value = ($plainPart + $italicPart)

$italicStart = $plainPart.Length + 1
$italicLen = $italicPart.Length
$ws.Range("C82").Characters($italicStart, $italicLen).Font.Italic = $true

# D82 stays blank (inherits the column's default wrap-text style automatically)

# E82: French translation
$ws.Range("E82").Value = "Suspension buvable/rectale"

# Select the newly-edited cell, matching the author's saved cursor position
$ws.Range("E84").Select()

# The _FilterDatabase named range covers column E and must grow by one row
# to keep including the newly inserted row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$E`$1:`$E`$152"
    }
}
